$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# Insert a new row at 44, shifting existing rows 44+ down to 45+
$ws.Rows.Item(44).Insert()

# Split the combined wind row (old row 43: E[_]W* / wind) into
# offshore (row 43) and onshore (new row 44) entries
$ws.Range("D43").Value = "windoff"

$ws.Range("C44").Value = "E[_]WON*"
$ws.Range("D44").Value = "windon"
$ws.Range("E44").Value = "IN"

$ws.Range("C43").Value = "E[_]WOF*"

# Rename ELC_won* -> ELC_wo* (row 41, ElcAgg_Wind aggregation pattern)
$ws.Range("D41").Value = "ELC_wo*"

# Update sheet view to match the saved workbook state
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D42").Select()
